$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.073631167411804
$ws.Range("B1").Value = 1.773496150970459
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.873334407806396
$ws.Range("E1").Value = 1.156263589859009
